# Auto-generated edit script: refresh Leve profitability market-data columns (H-N)
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 724.41174
$ws.Range("I33").Value = 520.4231
$ws.Range("J33").Value = 1387.375
$ws.Range("K33").Value = 520.4231
$ws.Range("L33").Value = 1387.375
$ws.Range("M33").Value = -291.4231
$ws.Range("N33").Value = -1845.375

# Row 40
$ws.Range("H40").Value = 1099.7778
$ws.Range("I40").Value = 900
$ws.Range("J40").Value = 1259.6
$ws.Range("K40").Value = 900
$ws.Range("L40").Value = 1259.6
$ws.Range("M40").Value = -725
$ws.Range("N40").Value = -1609.6

# Row 129
$ws.Range("H129").Value = 1123
$ws.Range("J129").Value = 1320.2285
$ws.Range("L129").Value = 3960.6855
$ws.Range("N129").Value = -13960.6855

# Row 138
$ws.Range("H138").Value = 2097.1477
$ws.Range("J138").Value = 1991
$ws.Range("L138").Value = 5973
$ws.Range("N138").Value = -16253

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 606329.9399999999
$ws.Range("I32").Value = 739675.0600000001
$ws.Range("J32").Value = 19611.4
$ws.Range("K32").Value = 739675.0600000001
$ws.Range("L32").Value = 19611.4
$ws.Range("M32").Value = -739388.0600000001
$ws.Range("N32").Value = -20185.4

# Row 45
$ws.Range("H45").Value = 3113.3333
$ws.Range("I45").Value = 2882.4
$ws.Range("K45").Value = 2882.4
$ws.Range("M45").Value = -2505.4

# Row 122
$ws.Range("H122").Value = 2008.9714
$ws.Range("I122").Value = 1921.9259
$ws.Range("J122").Value = 2302.75
$ws.Range("K122").Value = 5765.7777
$ws.Range("L122").Value = 6908.25
$ws.Range("M122").Value = -3315.7777
$ws.Range("N122").Value = -11808.25

# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

# Row 132
$ws.Range("H132").Value = 5057.4346
$ws.Range("I132").Value = 5693
$ws.Range("J132").Value = 4474.8335
$ws.Range("K132").Value = 17079
$ws.Range("L132").Value = 13424.5005
$ws.Range("M132").Value = -14549
$ws.Range("N132").Value = -18484.5005

$ws = $wb.Worksheets.Item("BSM")
# Row 132
$ws.Range("H132").Value = 44026.668
$ws.Range("J132").Value = 44026.668
$ws.Range("L132").Value = 44026.668
$ws.Range("N132").Value = -54146.668

# Row 134
$ws.Range("H134").Value = 3534.6667
$ws.Range("I134").Value = 3476.5
$ws.Range("K134").Value = 10429.5
$ws.Range("M134").Value = -7894.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 7295.4165
$ws.Range("I31").Value = 1660.0454
$ws.Range("J31").Value = 16151
$ws.Range("K31").Value = 1660.0454
$ws.Range("L31").Value = 16151
$ws.Range("M31").Value = -1365.0454
$ws.Range("N31").Value = -16741

# Row 34
$ws.Range("H34").Value = 7295.4165
$ws.Range("I34").Value = 1660.0454
$ws.Range("J34").Value = 16151
$ws.Range("K34").Value = 1660.0454
$ws.Range("L34").Value = 16151
$ws.Range("M34").Value = -1458.0454
$ws.Range("N34").Value = -16555

# Row 51
$ws.Range("H51").Value = 12248.75
$ws.Range("J51").Value = 12248.75
$ws.Range("L51").Value = 12248.75
$ws.Range("N51").Value = -13720.75

# Row 59
$ws.Range("H59").Value = 25832.666
$ws.Range("J59").Value = 25832.666
$ws.Range("L59").Value = 25832.666
$ws.Range("N59").Value = -28122.666

# Row 60
$ws.Range("H60").Value = 10501.5
$ws.Range("J60").Value = 10501.5
$ws.Range("L60").Value = 10501.5
$ws.Range("N60").Value = -11523.5

# Row 61
$ws.Range("H61").Value = 12248.75
$ws.Range("J61").Value = 12248.75
$ws.Range("L61").Value = 12248.75
$ws.Range("N61").Value = -12944.75

# Row 68
$ws.Range("H68").Value = 21813.715
$ws.Range("I68").Value = 10200
$ws.Range("J68").Value = 23749.334
$ws.Range("K68").Value = 10200
$ws.Range("L68").Value = 23749.334
$ws.Range("M68").Value = -9451
$ws.Range("N68").Value = -25247.334

# Row 71
$ws.Range("H71").Value = 21813.715
$ws.Range("I71").Value = 10200
$ws.Range("J71").Value = 23749.334
$ws.Range("K71").Value = 30600
$ws.Range("L71").Value = 71248.00199999999
$ws.Range("M71").Value = -26856
$ws.Range("N71").Value = -78736.00199999999

# Row 74
$ws.Range("H74").Value = 25110.666
$ws.Range("I74").Value = 10000
$ws.Range("J74").Value = 29428
$ws.Range("K74").Value = 10000
$ws.Range("L74").Value = 29428
$ws.Range("M74").Value = -9126
$ws.Range("N74").Value = -31176

# Row 77
$ws.Range("H77").Value = 25110.666
$ws.Range("I77").Value = 10000
$ws.Range("J77").Value = 29428
$ws.Range("K77").Value = 30000
$ws.Range("L77").Value = 88284
$ws.Range("M77").Value = -25632
$ws.Range("N77").Value = -97020

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 11005
$ws.Range("I3").Value = 2030
$ws.Range("K3").Value = 6090
$ws.Range("M3").Value = -5978

# Row 23
$ws.Range("H23").Value = 111111240
$ws.Range("I23").Value = 137.25
$ws.Range("J23").Value = 200000130
$ws.Range("K23").Value = 411.75
$ws.Range("L23").Value = 600000390
$ws.Range("M23").Value = -176.75
$ws.Range("N23").Value = -600000860

# Row 39
$ws.Range("H39").Value = 1432.5641
$ws.Range("J39").Value = 1432.5641
$ws.Range("L39").Value = 4297.692300000001
$ws.Range("N39").Value = -4885.692300000001

# Row 110
$ws.Range("H110").Value = 14214.7
$ws.Range("I110").Value = 3013.5
$ws.Range("J110").Value = 17015
$ws.Range("K110").Value = 9040.5
$ws.Range("L110").Value = 51045
$ws.Range("M110").Value = -4950.5
$ws.Range("N110").Value = -59225

# Row 113
$ws.Range("H113").Value = 836.2727
$ws.Range("I113").Value = 458.25
$ws.Range("K113").Value = 1374.75
$ws.Range("M113").Value = 795.25

# Row 134
$ws.Range("H134").Value = 6204
$ws.Range("I134").Value = 5412
$ws.Range("J134").Value = 6600
$ws.Range("K134").Value = 16236
$ws.Range("L134").Value = 19800
$ws.Range("M134").Value = -11166
$ws.Range("N134").Value = -29940

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5463.5854
$ws.Range("I70").Value = 5469.231
$ws.Range("J70").Value = 5453.8
$ws.Range("K70").Value = 5469.231
$ws.Range("L70").Value = 5453.8
$ws.Range("M70").Value = -5199.231
$ws.Range("N70").Value = -5993.8

# Row 73
$ws.Range("H73").Value = 5463.5854
$ws.Range("I73").Value = 5469.231
$ws.Range("J73").Value = 5453.8
$ws.Range("K73").Value = 5469.231
$ws.Range("L73").Value = 5453.8
$ws.Range("M73").Value = -4533.231
$ws.Range("N73").Value = -7325.8

# Row 98
$ws.Range("H98").Value = 27500
$ws.Range("J98").Value = 27500
$ws.Range("L98").Value = 27500
$ws.Range("N98").Value = -33490

# Row 122
$ws.Range("H122").Value = 4659.722
$ws.Range("I122").Value = 1825
$ws.Range("J122").Value = 5469.643
$ws.Range("K122").Value = 5475
$ws.Range("L122").Value = 16408.929
$ws.Range("M122").Value = -3025
$ws.Range("N122").Value = -21308.929

# Row 126
$ws.Range("H126").Value = 2953.9
$ws.Range("I126").Value = 1378
$ws.Range("J126").Value = 4004.5
$ws.Range("K126").Value = 4134
$ws.Range("L126").Value = 12013.5
$ws.Range("M126").Value = -1664
$ws.Range("N126").Value = -16953.5

$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3167.5
$ws.Range("J122").Value = 3401
$ws.Range("L122").Value = 10203
$ws.Range("N122").Value = -15103

# Row 132
$ws.Range("H132").Value = 3385.1924
$ws.Range("I132").Value = 2665.9412
$ws.Range("J132").Value = 4743.778
$ws.Range("K132").Value = 7997.823600000001
$ws.Range("L132").Value = 14231.334
$ws.Range("M132").Value = -5467.823600000001
$ws.Range("N132").Value = -19291.334

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 756.2857
$ws.Range("I107").Value = 732.3333
$ws.Range("K107").Value = 2196.9999
$ws.Range("M107").Value = -276.9998999999998

# Row 122
$ws.Range("H122").Value = 3572.4062
$ws.Range("I122").Value = 2081.6365
$ws.Range("J122").Value = 4353.2856
$ws.Range("K122").Value = 6244.9095
$ws.Range("L122").Value = 13059.8568
$ws.Range("M122").Value = -3794.9095
$ws.Range("N122").Value = -17959.8568

# Row 123
$ws.Range("H123").Value = 10390
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()
